$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..21) {
    $ws.Range("E$r").Value = "NA"
}
